$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 80, pushing the existing rows 80-108 down to 81-109.
$ws.Rows.Item(80).Insert()

# Populate the newly inserted row 80 with the new record.
$ws.Cells.Item(80, 1).Value  = 7
$ws.Cells.Item(80, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(80, 3).Value  = "Ñuble"
$ws.Cells.Item(80, 4).Value  = 45141
$ws.Cells.Item(80, 5).Value  = 16
$ws.Cells.Item(80, 6).Value  = 100112044
$ws.Cells.Item(80, 7).Value  = "Perejil"
$ws.Cells.Item(80, 8).Value  = "Sin especificar"
$ws.Cells.Item(80, 9).Value  = "Primera"
$ws.Cells.Item(80, 10).Value = 60
$ws.Cells.Item(80, 11).Value = 1500
$ws.Cells.Item(80, 12).Value = 1500
$ws.Cells.Item(80, 13).Value = 1500
$ws.Cells.Item(80, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(80, 15).Value = "Región de Ñuble"
$ws.Cells.Item(80, 16).Value = 1500
$ws.Cells.Item(80, 17).Value = 1
$ws.Cells.Item(80, 18).Value = "Hortaliza"
